# Auto-generated Excel COM-interop script
# Applies updated currentAveragePrice / LevePrice* / LeveProfit* values
# across multiple sheets, as captured by the commit's xml diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 67.85714
$ws.Range("I11").Value = 67.85714
$ws.Range("K11").Value = 67.85714
$ws.Range("M11").Value = 72.14286

# Row 38
$ws.Range("H38").Value = 336.25
$ws.Range("I38").Value = 148.33333
$ws.Range("K38").Value = 444.99999
$ws.Range("M38").Value = -72.99998999999997

# Row 40
$ws.Range("H40").Value = 1550
$ws.Range("I40").Value = 1650
$ws.Range("J40").Value = 1350
$ws.Range("K40").Value = 1650
$ws.Range("L40").Value = 1350
$ws.Range("M40").Value = -1475
$ws.Range("N40").Value = -1700

# Row 64
$ws.Range("H64").Value = 3017.9167
$ws.Range("I64").Value = 3100
$ws.Range("J64").Value = 2976.875
$ws.Range("K64").Value = 3100
$ws.Range("L64").Value = 2976.875
$ws.Range("M64").Value = -2852
$ws.Range("N64").Value = -3472.875

# Row 67
$ws.Range("H67").Value = 3017.9167
$ws.Range("I67").Value = 3100
$ws.Range("J67").Value = 2976.875
$ws.Range("K67").Value = 3100
$ws.Range("L67").Value = 2976.875
$ws.Range("M67").Value = -2242
$ws.Range("N67").Value = -4692.875

# Row 69
$ws.Range("H69").Value = 3830.9285
$ws.Range("I69").Value = 3013
$ws.Range("J69").Value = 3893.8462
$ws.Range("K69").Value = 9039
$ws.Range("L69").Value = 11681.5386
$ws.Range("N69").Value = -13429.5386
$ws.Range("M69").Value = -8165

# Row 72
$ws.Range("H72").Value = 3830.9285
$ws.Range("I72").Value = 3013
$ws.Range("J72").Value = 3893.8462
$ws.Range("K72").Value = 27117
$ws.Range("L72").Value = 35044.6158
$ws.Range("N72").Value = -43780.6158
$ws.Range("M72").Value = -22749

# Row 112
$ws.Range("H112").Value = 5750.9375
$ws.Range("J112").Value = 6307.931
$ws.Range("L112").Value = 18923.793
$ws.Range("N112").Value = -21139.793

# Row 135
$ws.Range("H135").Value = 1116.5555
$ws.Range("I135").Value = 582.5769
$ws.Range("J135").Value = 15000
$ws.Range("K135").Value = 5243.1921
$ws.Range("L135").Value = 135000
$ws.Range("M135").Value = -2708.1921
$ws.Range("N135").Value = -140070

# Row 136
$ws.Range("H136").Value = 35000
$ws.Range("J136").Value = 35000
$ws.Range("L136").Value = 35000
$ws.Range("N136").Value = -45200

# Row 137
$ws.Range("H137").Value = 2933.5557
$ws.Range("I137").Value = 4667.3335
$ws.Range("J137").Value = 2066.6667
$ws.Range("K137").Value = 14002.0005
$ws.Range("L137").Value = 6200.000100000001
$ws.Range("M137").Value = -11452.0005
$ws.Range("N137").Value = -11300.0001

# Row 138
$ws.Range("H138").Value = 2589.359
$ws.Range("I138").Value = 2668.6667
$ws.Range("J138").Value = 2565.5667
$ws.Range("K138").Value = 8006.000100000001
$ws.Range("L138").Value = 7696.7001
$ws.Range("M138").Value = -2866.000100000001
$ws.Range("N138").Value = -17976.7001

# Row 139
$ws.Range("H139").Value = 44923.332
$ws.Range("J139").Value = 44923.332
$ws.Range("L139").Value = 44923.332
$ws.Range("N139").Value = -55203.332

# Row 140
$ws.Range("H140").Value = 77000
$ws.Range("J140").Value = 77000
$ws.Range("L140").Value = 77000
$ws.Range("N140").Value = -87360

# Row 141
$ws.Range("H141").Value = 2288.0308
$ws.Range("I141").Value = 630.84
$ws.Range("J141").Value = 7812
$ws.Range("K141").Value = 1892.52
$ws.Range("L141").Value = 23436
$ws.Range("M141").Value = 3287.48
$ws.Range("N141").Value = -33796

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 847456.5600000001
$ws.Range("I32").Value = 1019778.2
$ws.Range("J32").Value = 14568.833
$ws.Range("K32").Value = 1019778.2
$ws.Range("L32").Value = 14568.833
$ws.Range("M32").Value = -1019491.2
$ws.Range("N32").Value = -15142.833

# Row 132
$ws.Range("H132").Value = 2293.4834
$ws.Range("I132").Value = 1834.6351
$ws.Range("J132").Value = 4290.8237
$ws.Range("K132").Value = 5503.9053
$ws.Range("L132").Value = 12872.4711
$ws.Range("M132").Value = -2973.9053
$ws.Range("N132").Value = -17932.4711

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 3805.25
$ws.Range("I134").Value = 3869.5334
$ws.Range("J134").Value = 3731.077
$ws.Range("K134").Value = 11608.6002
$ws.Range("L134").Value = 11193.231
$ws.Range("M134").Value = -9073.600199999999
$ws.Range("N134").Value = -16263.231

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4301.2
$ws.Range("I31").Value = 1109.6428
$ws.Range("J31").Value = 7610.963
$ws.Range("K31").Value = 1109.6428
$ws.Range("L31").Value = 7610.963
$ws.Range("M31").Value = -814.6428000000001
$ws.Range("N31").Value = -8200.963

# Row 34
$ws.Range("H34").Value = 4301.2
$ws.Range("I34").Value = 1109.6428
$ws.Range("J34").Value = 7610.963
$ws.Range("K34").Value = 1109.6428
$ws.Range("L34").Value = 7610.963
$ws.Range("M34").Value = -907.6428000000001
$ws.Range("N34").Value = -8014.963

# Row 134
$ws.Range("H134").Value = 2787.7378
$ws.Range("I134").Value = 3177.3618
$ws.Range("J134").Value = 1479.7142
$ws.Range("K134").Value = 9532.0854
$ws.Range("L134").Value = 4439.142599999999
$ws.Range("M134").Value = -6997.0854
$ws.Range("N134").Value = -9509.142599999999

$ws = $wb.Worksheets.Item("CUL")
# Row 44
$ws.Range("H44").Value = 324.2857
$ws.Range("J44").Value = 400
$ws.Range("L44").Value = 1200
$ws.Range("N44").Value = -1996

# Row 131
$ws.Range("H131").Value = 2859.4844
$ws.Range("I131").Value = 510.1111
$ws.Range("J131").Value = 3243.9272
$ws.Range("K131").Value = 1530.3333
$ws.Range("L131").Value = 9731.7816
$ws.Range("M131").Value = 3509.6667
$ws.Range("N131").Value = -19811.7816

$ws = $wb.Worksheets.Item("GSM")
# Row 12
$ws.Range("H12").Value = 52102.4
$ws.Range("J12").Value = 65003
$ws.Range("L12").Value = 65003
$ws.Range("N12").Value = -65283

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 2072.4583
$ws.Range("I132").Value = 1932.9231
$ws.Range("J132").Value = 2435.25
$ws.Range("K132").Value = 5798.7693
$ws.Range("L132").Value = 7305.75
$ws.Range("M132").Value = -3268.7693
$ws.Range("N132").Value = -12365.75

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 79500
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 79500
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 79500
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -80748

# Row 65
$ws.Range("H65").Value = 79500
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 79500
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 397500
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -403740

# Row 132
$ws.Range("H132").Value = 2874597
$ws.Range("I132").Value = 959.3333
$ws.Range("J132").Value = 10417896
$ws.Range("K132").Value = 2877.9999
$ws.Range("L132").Value = 31253688
$ws.Range("M132").Value = -347.9998999999998
$ws.Range("N132").Value = -31258748
